# Update Runmode flags on the "Test Cases" sheet and refresh the IR/config
# codes on the incident-report / quality-dashboard sheets.

$wb = $excel.ActiveWorkbook

# --- "Test Cases" sheet: flip which test cases are flagged to run ---
$testCases = $wb.Worksheets.Item("Test Cases")

# TC_IncidentReport is now enabled for this run
$testCases.Range("C2").Value = "Yes"

# The sentinel/patient-complaint/quality-dashboard related cases are
# disabled again
$testCases.Range("C17").Value = "No"
$testCases.Range("C18").Value = "No"
$testCases.Range("C19").Value = "No"
$testCases.Range("C20").Value = "No"
$testCases.Range("C21").Value = "No"
$testCases.Range("C22").Value = "No"

# --- TC_IncidentReport sheet: bump the generated IR / config code ---
$incidentReport = $wb.Worksheets.Item("TC_IncidentReport")
$incidentReport.Range("G2").Value = "325 /2022 CONFIG CODE"

# --- TC_QualityDashboardAssign sheet: bump the IR code used for assign ---
$qdAssign = $wb.Worksheets.Item("TC_QualityDashboardAssign")
$qdAssign.Range("A2").Value = "326 /2022 CONFIG CODE"

# --- TC_QualityDashboardReAssign sheet: bump the IR code used for re-assign ---
$qdReAssign = $wb.Worksheets.Item("TC_QualityDashboardReAssign")
$qdReAssign.Range("A2").Value = "326 /2022 CONFIG CODE"
